$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# Sheet3: previously empty, now holds crop-tracking data for the ARB.
# Populate the new headers first so the new shared-string entries ("Crop
# Type", "Start Date", "End Date") are minted before the Sheet2/Sheet1
# edits below.
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "Crop Type"
$ws3.Range("B1").Value = "Start Date"
$ws3.Range("C1").Value = "End Date"
$ws3.Range("D1").Value = "ARB Name"
$ws3.Range("A1:D1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Sheet2: insert a new "Relationship" column (D) between "Education Level"
# (C) and "ARB Name" (old D, now E); update Education Level values to
# "Vocational"
# ---------------------------------------------------------------------------
$ws2.Columns("D").Insert()

$ws2.Range("D1").Value = "Relationship"
$ws2.Range("D2").Value = "Mother"
$ws2.Range("D3").Value = "Father"
$ws2.Range("D4").Value = "Wife"
$ws2.Range("D5").Value = "Brother"
$ws2.Range("D6").Value = "Sister"

# ---------------------------------------------------------------------------
# Sheet1: "Education Level" values change from "Tertiary - Undergraduate"
# to "College Level"
# ---------------------------------------------------------------------------
$ws1.Range("K2:K5").Value = "College Level"
$ws1.Columns("K").AutoFit()

# Sheet2 "Education Level" values change to "Vocational"
$ws2.Range("C2:C6").Value = "Vocational"

# ---------------------------------------------------------------------------
# Finish populating Sheet3's data row (crop type "Rice" is the very last new
# shared string the original workbook introduces).
# ---------------------------------------------------------------------------
$ws3.Range("B2").Value = 43152
$ws3.Range("C2").Value = 43187
$ws3.Range("D2").Value = "Rey Christian Lopez Gamboa"
$ws3.Range("A2").Value = "Rice"

# Match the date formatting already used elsewhere in the workbook (copy the
# number format from an existing "Member Since"/"Birthday" date cell instead
# of minting a fresh custom number format).
$ws1.Range("E2").Copy()
$ws3.Range("B2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping
# ---------------------------------------------------------------------------
$ws2.Range("D13").Select()
$ws3.Range("D13").Select()
$ws1.Range("D16").Select()
$ws1.Activate()
